# Append the latest scraped job-numbers snapshot as a new row at the
# bottom of the data table (job watcher data retrieval commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3

$ws.Cells.Item($row, 1).Value = 44316.7736188719   # date (2021-04-30 18:34:00 UTC)
$ws.Cells.Item($row, 2).Value = 71841              # tot
$ws.Cells.Item($row, 3).Value = 60264              # chde
$ws.Cells.Item($row, 4).Value = 3276               # bank
$ws.Cells.Item($row, 5).Value = 1985               # strat
$ws.Cells.Item($row, 6).Value = 1396               # strat_de
$ws.Cells.Item($row, 7).Value = 18793              # zh
$ws.Cells.Item($row, 8).Value = 1381               # bank_zh
$ws.Cells.Item($row, 9).Value = 792                # strat_zh
$ws.Cells.Item($row, 10).Value = 199               # strat_bank_zh
